$wb = $excel.ActiveWorkbook

# --- Sheet "compare_models": update column I (TT (Sec)) values ---
$wsCompare = $wb.Worksheets.Item("compare_models")

$compareUpdates = @{
    2  = 0.066
    3  = 0.058
    4  = 0.08
    5  = 0.038
    6  = 0.04
    7  = 0.02
    8  = 0.018
    9  = 0.018
    10 = 0.992
    11 = 0.508
    12 = 0.016
    14 = 0.026
    16 = 0.018
    17 = 0.026
    18 = 0.02
    19 = 0.02
}

foreach ($row in $compareUpdates.Keys) {
    $wsCompare.Range("I$row").Value = $compareUpdates[$row]
}

# --- Sheet "pred_final": update row 2 metrics (C2:H2) ---
$wsPredFinal = $wb.Worksheets.Item("pred_final")

$wsPredFinal.Range("C2").Value = 1.2839
$wsPredFinal.Range("D2").Value = 4.018
$wsPredFinal.Range("E2").Value = 2.0045
$wsPredFinal.Range("F2").Value = 0.994
$wsPredFinal.Range("G2").Value = 0.0337
$wsPredFinal.Range("H2").Value = 0.0211
